$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Type de mission" column (F) with its header and values
$ws.Range("F1").Value = "Type de mission"
$ws.Range("F2").Value = "AIR - GROUND"
$ws.Range("F3").Value = "AIR - AIR"
$ws.Range("F4").Value = "AIR - AIR"

# Set the width of the new column F (closest achievable value to 19.1640625
# given this engine's column-width quantization)
$ws.Columns("F").ColumnWidth = 18.33

# Update the active selection to match the edited workbook's state
$ws.Range("F7").Select()
